# "correct unit test for dao layer"
# The test fixture's header row used SCREAMING_CASE column-name labels
# (ID, USERNAME, PASSWORD, DISPLAY_NAME, REAL_NAME) that don't match the
# lowercase column names the dao layer actually expects/returns. Fix the
# header row on the "bg_user" sheet to use the correct lowercase names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bg_user")

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "display_name"
$ws.Range("E1").Value = "real_name"
